$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new row of data (row 96) continuing the existing series in A:C,
# copying the date-formatted style from the cell above (A95) onto A96.
$ws.Range("A95").Copy($ws.Range("A96"))

$ws.Range("A96").Value = 45597
$ws.Range("B96").Value = 0.193063140749257
$ws.Range("C96").Value = 0.0714225677663375
